$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2,290"; New = "2,416" },
    @{ Old = "6.2";   New = "6.6" },
    @{ Old = "3,333"; New = "3,368" },
    @{ Old = "11.1";  New = "11.2" },
    @{ Old = "753";   New = "786" },
    @{ Old = "18.4";  New = "19.2" },
    @{ Old = "1,556"; New = "1,567" },
    @{ Old = "521";   New = "562" },
    @{ Old = "1,653"; New = "1,771" },
    @{ Old = "336";   New = "360" },
    @{ Old = "7,011"; New = "7,205" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
